$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently ends with a thick-bottom-border "closer" row at row 5
# (A5/B5, empty). We need to:
#   - turn row 5 into a normal data row: "Ростислав Бердниченко" / "Development"
#   - add a new closer row 6: "Науменко Артем" / "Writing documentation"
# Easiest robust way: insert a new row above the current closer row, copying
# its formatting from the row above (a normal data row), then fill both rows.

$ws.Rows.Item(5).Insert()

# Copy the normal-row formatting (row 4) into the newly inserted row 5
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "Ростислав Бердниченко"
$ws.Range("B5").Value = "Development"

$ws.Range("A6").Value = "Науменко Артем"
$ws.Range("B6").Value = "Writing documentation"

# Column A widened (bestFit) to accommodate the longer name just added.
$ws.Columns.Item(1).ColumnWidth = 21.15

$ws.Range("I14").Select()
